# Apply the "glmm_model_selection_hab_season_lmb" update:
#  - drop the old "deviance" column (old column L) entirely, shifting the
#    old "df.residual" column (old column M) into its place (new column L)
#  - switch the fitted family from Gamma to lognormal everywhere
#  - append the ar1() random-effect terms to the top model's formula (row 2)
#  - refresh all of the re-fit statistics (sigma, logLik, AIC, delta_AIC,
#    AIC_weight, BIC, df.residual, nobs) to their new values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "deviance" column; this shifts the old "df.residual"
# column (M) left into column L, header included, matching the diff.
$ws.Columns("L").Delete()

# --- Row 2 (model "m") ---
$ws.Range("A2").Value = "lognormal"
$ws.Range("C2").Value = "mean_accel ~ habitat_type * season + (1 | animal_id) + ar1(season + 0 | animal_id) + ar1(habitat_type + 0 | animal_id)"
$ws.Range("F2").Value = 0.367886140934635
$ws.Range("G2").Value = 3453.37736771149
$ws.Range("H2").Value = -6858.75473542298
$ws.Range("K2").Value = -6673.2212542058
$ws.Range("L2").Value = 16800

# --- Row 3 (model "m2") ---
$ws.Range("A3").Value = "lognormal"
$ws.Range("E3").Value = 16824
$ws.Range("F3").Value = 0.373938160111812
$ws.Range("G3").Value = 3280.52269517336
$ws.Range("H3").Value = -6549.04539034673
$ws.Range("I3").Value = 309.709345076252
$ws.Range("J3").Value = 0.0000000000000000000000000000000000000000000000000000000000000000000559075189607128
$ws.Range("K3").Value = -6502.66202004243
$ws.Range("L3").Value = 16818

# --- Row 4 (model "m1") ---
$ws.Range("A4").Value = "lognormal"
$ws.Range("F4").Value = 0.375750086752719
$ws.Range("G4").Value = 2140.6878653588
$ws.Range("H4").Value = -4267.3757307176
$ws.Range("I4").Value = 2591.37900470538
$ws.Range("K4").Value = -4213.26179869592
$ws.Range("L4").Value = 16817

$wb.Save()
